# Weekly price-sheet update: a new price record for "Damasco" was reported
# for the Provincia de San Felipe de Aconcagua market and needs to be
# inserted as the new row 4 (directly below the two most-recent records),
# pushing every existing record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 4 - this shifts rows 4..56 down to 5..57
# and keeps the date-number formatting from the row above for column D.
$ws.Rows(4).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(4, 1).Value  = 10
$ws.Cells.Item(4, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(4, 3).Value  = "La Araucanía"
$ws.Cells.Item(4, 4).Value  = 44524
$ws.Cells.Item(4, 5).Value  = 9
$ws.Cells.Item(4, 6).Value  = "Fruta"
$ws.Cells.Item(4, 7).Value  = 100103
$ws.Cells.Item(4, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(4, 9).Value  = 100103003
$ws.Cells.Item(4, 10).Value = "Damasco"
$ws.Cells.Item(4, 11).Value = "Castle Brite"
$ws.Cells.Item(4, 12).Value = "Primera"
$ws.Cells.Item(4, 13).Value = 145
$ws.Cells.Item(4, 14).Value = 18000
$ws.Cells.Item(4, 15).Value = 20000
$ws.Cells.Item(4, 16).Value = 19103
$ws.Cells.Item(4, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(4, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(4, 19).Value = 2729
$ws.Cells.Item(4, 20).Value = 7
